$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{Row=2; D='67.439.92'; E='  +1.18%  '},
    @{Row=3; D='3.941.50'; E='  +4.32%  '},
    @{Row=4; E='  +0.07%  '},
    @{Row=5; D='470.69'; E='  +7.68%  '},
    @{Row=6; D='147.49'; E='  +3.96%  '},
    @{Row=7; E='  +0.63%  '},
    @{Row=8; E='  -0.03%  '},
    @{Row=9; D='0.730'; E='  -0.87%  '},
    @{Row=10; D='0.168'; E='  +10.37%  '},
    @{Row=11; D='0.0000350'; E='  +10.82%  '},
    @{Row=12; E='  +0.86%  '},
    @{Row=13; D='4.566.34'; E='  +4.15%  '},
    @{Row=14; D='10.36'; E='  -0.75%  '},
    @{Row=15; D='15.07'; E='  +1.87%  '},
    @{Row=16; D='3.933.35'; E='  +3.76%  '},
    @{Row=17; E='  +0.40%  '},
    @{Row=18; D='19.87'; E='  +0.07%  '},
    @{Row=19; E='  +1.85%  '},
    @{Row=20; D='67.565.38'; E='  +1.24%  '},
    @{Row=21; D='434.15'; E='  +4.18%  '},
    @{Row=22; D='3.38'; E='  +3.58%  '},
    @{Row=23; D='14.42'; E='  -0.58%  '},
    @{Row=24; D='87.48'; E='  +1.56%  '},
    @{Row=25; D='3.61'; E='  +6.38%  '},
    @{Row=26; D='38.55'; E='  +3.78%  '},
    @{Row=27; B='RenderToken'; C='https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'; D='9.91'; E='  +3.85%  '},
    @{Row=28; B='Filecoin'; C='https://coinranking.com/coin/ymQub4fuB+filecoin-fil'; D='10.12'; E='  +3.48%  '},
    @{Row=29; B='LEO'; C='https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'; D='5.67'; E='  +0.55%  '},
    @{Row=30; B='Bittensor'; C='https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'; D='722.62'; E='  +0.24%  '},
    @{Row=31; B='Hedera'; C='https://coinranking.com/coin/jad286TjB+hedera-hbar'; D='0.132'; E='  -1.30%  '},
    @{Row=32; B='Cosmos'; C='https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'; D='13.41'; E='  -2.71%  '},
    @{Row=33; B='Toncoin'; C='https://coinranking.com/coin/67YlI0K1b+toncoin-ton'; D='2.81'; E='  +2.33%  '},
    @{Row=34; B='InjectiveProtocol'; C='https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'; D='42.27'; E='  -2.84%  '},
    @{Row=35; D='57.94'; E='  +2.59%  '},
    @{Row=36; B='PEPE'; C='https://coinranking.com/coin/03WI8NQPF+pepe-pepe'; D='0.0₃0825'; E='  +22.35%  '},
    @{Row=37; B='Kaspa'; C='https://coinranking.com/coin/V8GxkwWow+kaspa-kas'; D='0.151'; E='  -2.92%  '},
    @{Row=38; B='Dai'; C='https://coinranking.com/coin/MoTuySvg7+dai-dai'; D='0.999'; E='  -0.10%  '},
    @{Row=39; B='NEARProtocol'; C='https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'; D='5.34'; E='  -4.69%  '},
    @{Row=40; B='VeChain'; C='https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'; D='0.0475'; E='  -0.43%  '},
    @{Row=41; B='ThetaToken'; C='https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta'; D='3.04'; E='  +4.19%  '},
    @{Row=42; B='Stellar'; C='https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'; D='0.142'; E='  +0.66%  '},
    @{Row=43; E='  -0.19%  '},
    @{Row=44; D='0.336'; E='  +3.09%  '},
    @{Row=45; B='LidoDAOToken'; C='https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'; D='3.49'; E='  +5.45%  '},
    @{Row=46; B='WEMIXToken'; C='https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'; D='2.82'; E='  +6.93%  '},
    @{Row=47; B='ARBITRUM'; C='https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'; D='2.21'; E='  +6.30%  '},
    @{Row=48; B='Fetch.AI'; C='https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'; D='2.54'; E='  -3.56%  '},
    @{Row=49; D='3.25'; E='  -2.92%  '},
    @{Row=50; B='Monero'; C='https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'; D='148.15'; E='  +3.93%  '},
    @{Row=51; B='Stacks'; C='https://coinranking.com/coin/mMPrMcB7+stacks-stx'; D='2.87'; E='  +1.47%  '}
)

foreach ($u in $updates) {
    if ($u.ContainsKey("B")) { $ws.Cells.Item($u.Row, 2).Value = $u.B }
    if ($u.ContainsKey("C")) { $ws.Cells.Item($u.Row, 3).Value = $u.C }
    if ($u.ContainsKey("D")) { $ws.Cells.Item($u.Row, 4).Value = $u.D }
    if ($u.ContainsKey("E")) { $ws.Cells.Item($u.Row, 5).Value = $u.E }
}
